$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 944
$ws.Range("I16").Value = 944
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 944
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -714
$ws.Range("N16").ClearContents()

# Row 132
$ws.Range("H132").Value = 1495.037
$ws.Range("I132").Value = 1495.037
$ws.Range("K132").Value = 4485.111
$ws.Range("M132").Value = -1955.111


$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3742.6667
$ws.Range("I45").Value = 3514.3333
$ws.Range("K45").Value = 3514.3333
$ws.Range("M45").Value = -3137.3333

# Row 61
$ws.Range("H61").Value = 1871.5385
$ws.Range("I61").Value = 1450.6666
$ws.Range("J61").Value = 2445.4546
$ws.Range("K61").Value = 1450.6666
$ws.Range("L61").Value = 2445.4546
$ws.Range("M61").Value = -1238.6666
$ws.Range("N61").Value = -2869.4546

# Row 97
$ws.Range("H97").Value = 327
$ws.Range("I97").Value = 327
$ws.Range("K97").Value = 327
$ws.Range("M97").Value = 169

# Row 132
$ws.Range("H132").Value = 2563.3076
$ws.Range("I132").Value = 1826.4445
$ws.Range("K132").Value = 5479.333500000001
$ws.Range("M132").Value = -2949.333500000001

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 135
$ws.Range("H135").Value = 59999.5
$ws.Range("J135").Value = 59999.5
$ws.Range("L135").Value = 59999.5
$ws.Range("N135").Value = -70139.5

# Row 136
$ws.Range("H136").Value = 1871.5385
$ws.Range("I136").Value = 1450.6666
$ws.Range("J136").Value = 2445.4546
$ws.Range("K136").Value = 4351.9998
$ws.Range("L136").Value = 7336.3638
$ws.Range("M136").Value = -1801.9998
$ws.Range("N136").Value = -12436.3638

# Row 138
$ws.Range("H138").Value = 99899.5
$ws.Range("J138").Value = 99899.5
$ws.Range("L138").Value = 99899.5
$ws.Range("N138").Value = -110179.5


$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 5347.4
$ws.Range("I11").Value = 169
$ws.Range("J11").Value = 8799.666999999999
$ws.Range("K11").Value = 169
$ws.Range("L11").Value = 8799.666999999999
$ws.Range("M11").Value = -29
$ws.Range("N11").Value = -9079.666999999999

# Row 81
$ws.Range("H81").Value = 74970
$ws.Range("J81").Value = 74970
$ws.Range("L81").Value = 74970
$ws.Range("N81").Value = -77092

# Row 84
$ws.Range("H84").Value = 74970
$ws.Range("J84").Value = 74970
$ws.Range("L84").Value = 224910
$ws.Range("N84").Value = -235518

# Row 94
$ws.Range("H94").Value = 506.8889
$ws.Range("I94").Value = 395
$ws.Range("K94").Value = 395
$ws.Range("M94").Value = 56


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("I31").Value = 2190.0908
$ws.Range("J31").Value = 2998
$ws.Range("K31").Value = 2190.0908
$ws.Range("L31").Value = 2998
$ws.Range("M31").Value = -1895.0908
$ws.Range("N31").Value = -3588

# Row 34
$ws.Range("I34").Value = 2190.0908
$ws.Range("J34").Value = 2998
$ws.Range("K34").Value = 2190.0908
$ws.Range("L34").Value = 2998
$ws.Range("M34").Value = -1988.0908
$ws.Range("N34").Value = -3402

# Row 122
$ws.Range("H122").Value = 2753
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 4461.625
$ws.Range("I18").Value = 700
$ws.Range("J18").Value = 4999
$ws.Range("K18").Value = 2100
$ws.Range("L18").Value = 14997
$ws.Range("M18").Value = -1931
$ws.Range("N18").Value = -15335

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 63
$ws.Range("H63").Value = 500
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 500
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 127
$ws.Range("H127").Value = 1000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

# Row 130
$ws.Range("H130").Value = 7333.3335
$ws.Range("J130").Value = 10000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

# Row 131
$ws.Range("H131").Value = 2007.5834
$ws.Range("J131").Value = 2063.182
$ws.Range("L131").Value = 6189.545999999999
$ws.Range("N131").Value = -16269.546

# Row 137
$ws.Range("H137").Value = 4555.5835
$ws.Range("J137").Value = 5266.8
$ws.Range("L137").Value = 15800.4
$ws.Range("N137").Value = -26000.4

# Row 138
$ws.Range("H138").Value = 5524.5
$ws.Range("I138").Value = 5524.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 16573.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -11433.5
$ws.Range("N138").ClearContents()

# Row 139
$ws.Range("H139").Value = 1496.2858
$ws.Range("I139").Value = 1496.2858
$ws.Range("K139").Value = 4488.857400000001
$ws.Range("M139").Value = 651.1425999999992

# Row 140
$ws.Range("H140").Value = 1738
$ws.Range("I140").Value = 1738
$ws.Range("K140").Value = 5214
$ws.Range("M140").Value = -34


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -13900

# Row 139
$ws.Range("H139").Value = 75685.625
$ws.Range("J139").Value = 75685.625
$ws.Range("L139").Value = 75685.625
$ws.Range("N139").Value = -85965.625


$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 416.85715
$ws.Range("J55").Value = 399
$ws.Range("L55").Value = 399
$ws.Range("N55").Value = -745

# Row 93
$ws.Range("H93").Value = 2120.95
$ws.Range("I93").Value = 2089.8823
$ws.Range("J93").Value = 2297
$ws.Range("K93").Value = 2089.8823
$ws.Range("L93").Value = 2297
$ws.Range("M93").Value = -841.8823000000002
$ws.Range("N93").Value = -4793

# Row 122
$ws.Range("H122").Value = 7000
$ws.Range("I122").Value = 7000
$ws.Range("K122").Value = 21000
$ws.Range("M122").Value = -18550


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3219.889
$ws.Range("I122").Value = 3155.1428
$ws.Range("J122").Value = 3446.5
$ws.Range("K122").Value = 9465.428400000001
$ws.Range("L122").Value = 10339.5
$ws.Range("M122").Value = -7015.428400000001
$ws.Range("N122").Value = -15239.5

# Row 135
$ws.Range("H135").Value = 71950
$ws.Range("J135").Value = 71950
$ws.Range("L135").Value = 71950
$ws.Range("N135").Value = -82090

# Row 136
$ws.Range("H136").Value = 2099.25
$ws.Range("I136").Value = 1699.25
$ws.Range("K136").Value = 5097.75
$ws.Range("M136").Value = -2547.75

